$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 193, shifting existing rows 193:207 down to 194:208.
$ws.Rows.Item(193).Insert()

# Populate the newly inserted row 193 with the new weekly record.
$ws.Cells.Item(193, 1).Value = 10
$ws.Cells.Item(193, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(193, 3).Value = "La Araucanía"
$ws.Cells.Item(193, 4).Value = 45166
$ws.Cells.Item(193, 5).Value = 9
$ws.Cells.Item(193, 6).Value = 100112035
$ws.Cells.Item(193, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(193, 8).Value = "Sin especificar"
$ws.Cells.Item(193, 9).Value = "Primera"
$ws.Cells.Item(193, 10).Value = 100
$ws.Cells.Item(193, 11).Value = 25000
$ws.Cells.Item(193, 12).Value = 25000
$ws.Cells.Item(193, 13).Value = 25000
$ws.Cells.Item(193, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(193, 15).Value = "Región Metropolitana"
$ws.Cells.Item(193, 16).Value = 1667
$ws.Cells.Item(193, 17).Value = 15
$ws.Cells.Item(193, 18).Value = "Hortaliza"
